$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (old row 6); rows 2-5 keep their data (values updated below)
$ws.Rows.Item(6).Delete()

# Update data values for rows 2-5 (34 columns: A..AH) to the new dataset
$ws.Cells.Item(2, 1).Value = 45146.50694444445
$ws.Cells.Item(2, 2).Value = 5.237
$ws.Cells.Item(2, 3).Value = 1.607
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 2.478
$ws.Cells.Item(2, 6).Value = 3.053
$ws.Cells.Item(2, 7).Value = 2.519
$ws.Cells.Item(2, 8).Value = 5.488
$ws.Cells.Item(2, 9).Value = 1.554
$ws.Cells.Item(2, 10).Value = 0.9409999999999999
$ws.Cells.Item(2, 11).Value = 4.022
$ws.Cells.Item(2, 12).Value = 1.069
$ws.Cells.Item(2, 13).Value = 0.9379999999999999
$ws.Cells.Item(2, 14).Value = 0.6929999999999999
$ws.Cells.Item(2, 15).Value = 0.87
$ws.Cells.Item(2, 16).Value = 2.749
$ws.Cells.Item(2, 17).Value = 1.106
$ws.Cells.Item(2, 18).Value = 0.51
$ws.Cells.Item(2, 19).Value = 0.063
$ws.Cells.Item(2, 20).Value = 20.094
$ws.Cells.Item(2, 21).Value = 4.803
$ws.Cells.Item(2, 22).Value = 2.372
$ws.Cells.Item(2, 23).Value = 3.891
$ws.Cells.Item(2, 24).Value = 1.034
$ws.Cells.Item(2, 25).Value = 0.249
$ws.Cells.Item(2, 26).Value = 1.98
$ws.Cells.Item(2, 27).Value = 1.154
$ws.Cells.Item(2, 28).Value = 0.674
$ws.Cells.Item(2, 29).Value = 0.9399999999999999
$ws.Cells.Item(2, 30).Value = 3.129
$ws.Cells.Item(2, 31).Value = 2.866
$ws.Cells.Item(2, 32).Value = 3.128
$ws.Cells.Item(2, 33).Value = 0.446
$ws.Cells.Item(2, 34).Value = 1.43
$ws.Cells.Item(3, 1).Value = 45146.51388888889
$ws.Cells.Item(3, 2).Value = 10.78
$ws.Cells.Item(3, 3).Value = 7.314
$ws.Cells.Item(3, 4).Value = 0.159
$ws.Cells.Item(3, 5).Value = 19.564
$ws.Cells.Item(3, 6).Value = 16.853
$ws.Cells.Item(3, 7).Value = 8.041
$ws.Cells.Item(3, 8).Value = 24.27
$ws.Cells.Item(3, 9).Value = 10.938
$ws.Cells.Item(3, 10).Value = 5.132
$ws.Cells.Item(3, 11).Value = 8.691000000000001
$ws.Cells.Item(3, 12).Value = 7.984
$ws.Cells.Item(3, 13).Value = 8.172000000000001
$ws.Cells.Item(3, 14).Value = 2.427
$ws.Cells.Item(3, 15).Value = 7.013
$ws.Cells.Item(3, 16).Value = 11.086
$ws.Cells.Item(3, 17).Value = 5.966
$ws.Cells.Item(3, 18).Value = 0.452
$ws.Cells.Item(3, 19).Value = 0.19
$ws.Cells.Item(3, 20).Value = 106.674
$ws.Cells.Item(3, 21).Value = 20.818
$ws.Cells.Item(3, 22).Value = 7.275
$ws.Cells.Item(3, 23).Value = 14.289
$ws.Cells.Item(3, 24).Value = 7.098
$ws.Cells.Item(3, 25).Value = 0.986
$ws.Cells.Item(3, 26).Value = 12.441
$ws.Cells.Item(3, 27).Value = 6.039
$ws.Cells.Item(3, 28).Value = 5.036
$ws.Cells.Item(3, 29).Value = 6.028
$ws.Cells.Item(3, 30).Value = 9.419
$ws.Cells.Item(3, 31).Value = 1.133
$ws.Cells.Item(3, 32).Value = 21.078
$ws.Cells.Item(3, 33).Value = 3.662
$ws.Cells.Item(3, 34).Value = 8.327999999999999
$ws.Cells.Item(4, 1).Value = 45146.52083333334
$ws.Cells.Item(4, 2).Value = 17.552
$ws.Cells.Item(4, 3).Value = 12.742
$ws.Cells.Item(4, 4).Value = 0.421
$ws.Cells.Item(4, 5).Value = 35.692
$ws.Cells.Item(4, 6).Value = 29.987
$ws.Cells.Item(4, 7).Value = 13.589
$ws.Cells.Item(4, 8).Value = 49.045
$ws.Cells.Item(4, 9).Value = 19.871
$ws.Cells.Item(4, 10).Value = 9.114000000000001
$ws.Cells.Item(4, 11).Value = 14.213
$ws.Cells.Item(4, 12).Value = 14.453
$ws.Cells.Item(4, 13).Value = 15.033
$ws.Cells.Item(4, 14).Value = 4.232
$ws.Cells.Item(4, 15).Value = 12.829
$ws.Cells.Item(4, 16).Value = 19.128
$ws.Cells.Item(4, 17).Value = 10.712
$ws.Cells.Item(4, 18).Value = 0.43
$ws.Cells.Item(4, 19).Value = 0.403
$ws.Cells.Item(4, 20).Value = 192.707
$ws.Cells.Item(4, 21).Value = 36.743
$ws.Cells.Item(4, 22).Value = 12.369
$ws.Cells.Item(4, 23).Value = 24.986
$ws.Cells.Item(4, 24).Value = 12.844
$ws.Cells.Item(4, 25).Value = 1.725
$ws.Cells.Item(4, 26).Value = 24.434
$ws.Cells.Item(4, 27).Value = 10.688
$ws.Cells.Item(4, 28).Value = 9.183
$ws.Cells.Item(4, 29).Value = 10.875
$ws.Cells.Item(4, 30).Value = 15.85
$ws.Cells.Item(4, 31).Value = 0.716
$ws.Cells.Item(4, 32).Value = 44.204
$ws.Cells.Item(4, 33).Value = 6.723
$ws.Cells.Item(4, 34).Value = 14.957
$ws.Cells.Item(5, 1).Value = 45146.52777777778
$ws.Cells.Item(5, 2).Value = 13.99
$ws.Cells.Item(5, 3).Value = 10.21
$ws.Cells.Item(5, 4).Value = 0.33
$ws.Cells.Item(5, 5).Value = 28.6
$ws.Cells.Item(5, 6).Value = 24.01
$ws.Cells.Item(5, 7).Value = 10.87
$ws.Cells.Item(5, 8).Value = 42.6
$ws.Cells.Item(5, 9).Value = 15.91
$ws.Cells.Item(5, 10).Value = 7.31
$ws.Cells.Item(5, 11).Value = 11.33
$ws.Cells.Item(5, 12).Value = 11.6
$ws.Cells.Item(5, 13).Value = 12.04
$ws.Cells.Item(5, 14).Value = 3.39
$ws.Cells.Item(5, 15).Value = 10.29
$ws.Cells.Item(5, 16).Value = 15.33
$ws.Cells.Item(5, 17).Value = 8.57
$ws.Cells.Item(5, 18).Value = 0.35
$ws.Cells.Item(5, 19).Value = 0.31
$ws.Cells.Item(5, 20).Value = 152.73
$ws.Cells.Item(5, 21).Value = 29.48
$ws.Cells.Item(5, 22).Value = 9.880000000000001
$ws.Cells.Item(5, 23).Value = 20.03
$ws.Cells.Item(5, 24).Value = 10.3
$ws.Cells.Item(5, 25).Value = 1.38
$ws.Cells.Item(5, 26).Value = 20.63
$ws.Cells.Item(5, 27).Value = 8.56
$ws.Cells.Item(5, 28).Value = 7.36
$ws.Cells.Item(5, 29).Value = 8.710000000000001
$ws.Cells.Item(5, 30).Value = 12.66
$ws.Cells.Item(5, 31).Value = 0.52
$ws.Cells.Item(5, 32).Value = 38.49
$ws.Cells.Item(5, 33).Value = 5.39
$ws.Cells.Item(5, 34).Value = 11.98

# Adjust column widths: C, Q, AA, AC widen from 7 to 8 chars; AD narrows from 8 to 7 chars
$ws.Columns.Item(3).ColumnWidth = 7.15
$ws.Columns.Item(17).ColumnWidth = 7.15
$ws.Columns.Item(27).ColumnWidth = 7.15
$ws.Columns.Item(29).ColumnWidth = 7.15
$ws.Columns.Item(30).ColumnWidth = 6.15
